$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the paragraph that holds the existing JetBrains / Azure SQL
# hyperlink (rId10) - this is the anchor point for the new content.
# ------------------------------------------------------------------
$anchorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*jetbrains.com/help/pycharm/azure-sql-database.html*") {
        $anchorPara = $p
        break
    }
}

# Give the paragraph mark of that paragraph the Hyperlink character
# style, matching the pattern used elsewhere in the document for
# paragraphs that are immediately followed by a blank "Hyperlink
# styled" paragraph.
$anchorPara.Range.ParagraphStyle
$markRange = $anchorPara.Range
$markRange.Collapse(0)
$markRange.MoveEnd(1, 1)
$markRange.Style = $d.Styles.Item("Hyperlink")

# ------------------------------------------------------------------
# Insert a new blank paragraph right after it, and give its paragraph
# mark the Hyperlink style too.
# ------------------------------------------------------------------
$insertionPoint = $anchorPara.Range
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()

$blankPara = $d.Paragraphs.Item($anchorPara.Index + 1)
$blankMark = $blankPara.Range
$blankMark.Collapse(0)
$blankMark.MoveEnd(1, 1)
$blankMark.Style = $d.Styles.Item("Hyperlink")

# ------------------------------------------------------------------
# Insert the new Azure App Service hyperlink paragraph.
# ------------------------------------------------------------------
$blankPara.Range.Collapse(0)
$blankPara.Range.InsertParagraphAfter()
$linkPara = $d.Paragraphs.Item($blankPara.Index + 1)
$linkRange = $linkPara.Range
$linkRange.Collapse(0)
$d.Hyperlinks.Add($linkRange, `
    "https://learn.microsoft.com/en-us/azure/app-service/quickstart-python?tabs=flask%2Cwindows%2Cazure-portal%2Cvscode-deploy%2Cdeploy-instructions-azportal%2Cterminal-bash%2Cdeploy-instructions-zip-azcli", `
    "", "", `
    "https://learn.microsoft.com/en-us/azure/app-service/quickstart-python?tabs=flask%2Cwindows%2Cazure-portal%2Cvscode-deploy%2Cdeploy-instructions-azportal%2Cterminal-bash%2Cdeploy-instructions-zip-azcli") | Out-Null

# ------------------------------------------------------------------
# Empty paragraph after the new hyperlink.
# ------------------------------------------------------------------
$linkPara.Range.Collapse(0)
$linkPara.Range.InsertParagraphAfter()
$emptyPara = $d.Paragraphs.Item($linkPara.Index + 1)

# ------------------------------------------------------------------
# Paragraph with the plain-text code "0QI7PAHYPAJU248J$".
# ------------------------------------------------------------------
$emptyPara.Range.Collapse(0)
$emptyPara.Range.InsertParagraphAfter()
$codePara = $d.Paragraphs.Item($emptyPara.Index + 1)
$codePara.Range.Text = "0QI7PAHYPAJU248J$"

# ------------------------------------------------------------------
# Paragraph with the "kmcoesdhhj" run, styled like pasted web text.
# ------------------------------------------------------------------
$codePara.Range.Collapse(0)
$codePara.Range.InsertParagraphAfter()
$stylePara = $d.Paragraphs.Item($codePara.Index + 1)
$stylePara.Range.Text = "kmcoesdhhj"
$stylePara.Range.Font.Name = "Segoe UI"
$stylePara.Range.Font.Color = 3289907
$stylePara.Range.Font.Size = 10
$stylePara.Range.Shading.BackgroundPatternColor = 16777215
